$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: intervention_type
$values = @{
    1 = 'intervention_type'
    2 = 'DRUG'
    3 = 'PROCEDURE'
    4 = 'DRUG'
    5 = 'PROCEDURE'
    6 = 'PROCEDURE'
    7 = 'DRUG'
    8 = 'DRUG'
    9 = 'PROCEDURE'
    10 = 'PROCEDURE'
    11 = 'DRUG'
    12 = 'BIOLOGICAL'
    13 = 'DRUG'
    14 = 'PROCEDURE'
    15 = 'OTHER'
    16 = 'DRUG'
    17 = 'OTHER'
    18 = 'DRUG'
    19 = 'BIOLOGICAL'
    20 = 'DRUG'
    21 = 'PROCEDURE'
    22 = 'DRUG'
    23 = 'DRUG'
    24 = 'DRUG'
    25 = 'BIOLOGICAL'
    26 = 'DRUG'
    27 = 'OTHER'
    28 = 'PROCEDURE'
    29 = 'PROCEDURE'
    30 = 'DRUG'
    31 = 'COMBINATION_PRODUCT'
    32 = 'DRUG'
    33 = 'OTHER'
    34 = 'OTHER'
    35 = 'OTHER'
    36 = 'OTHER'
    37 = 'RADIATION'
    38 = 'PROCEDURE'
    39 = 'DRUG'
    40 = 'OTHER'
    41 = 'OTHER'
    42 = 'DRUG'
    43 = 'OTHER'
    44 = 'OTHER'
    45 = 'DIAGNOSTIC_TEST'
    46 = 'DEVICE'
    47 = 'OTHER'
    48 = ''
    49 = 'PROCEDURE'
    50 = 'DRUG'
    51 = 'OTHER'
    52 = 'PROCEDURE'
    53 = 'PROCEDURE'
    54 = 'OTHER'
    55 = 'DIAGNOSTIC_TEST'
    56 = 'DRUG'
    57 = 'OTHER'
    58 = 'DRUG'
    59 = 'DRUG'
    60 = 'DEVICE'
    61 = 'OTHER'
    62 = 'DRUG'
    63 = 'OTHER'
    64 = 'DEVICE'
    65 = 'DEVICE'
    66 = 'DRUG'
    67 = 'DRUG'
    68 = 'DEVICE'
    69 = 'OTHER'
    70 = 'DIAGNOSTIC_TEST'
    71 = 'DIETARY_SUPPLEMENT'
    72 = 'GENETIC'
    73 = 'OTHER'
    74 = 'DRUG'
    75 = 'OTHER'
    76 = 'GENETIC'
    77 = 'DEVICE'
    78 = 'DIAGNOSTIC_TEST'
    79 = 'PROCEDURE'
    80 = 'PROCEDURE'
    81 = 'OTHER'
    82 = 'OTHER'
    83 = 'DEVICE'
    84 = ''
    85 = 'BIOLOGICAL'
    86 = 'DRUG'
}

foreach ($row in 1..86) {
    $v = $values[$row]
    $cell = $ws.Cells.Item($row, 11)
    $cell.Value = $v
}

# Match the header style (bold + thin border + center/top align) used by A1:J1
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Done"
